# Updated notebook, reran simulation
#
# Two new simulation rows ("Holden" and "Rizzie Spiral") are inserted right
# after the existing "Spiral5" row, pushing every category that used to
# follow it down by two rows (same category order, same labels - except
# "Thomas Hex" is renamed to "Matthies Hex"). The rerun also produced new
# numeric results for every category, and two more rows ("Michael-CCHex",
# "Michael-SNHex") are appended at the end to keep the full set.

function Set-RowValues($ws, $rangeAddr, $values) {
    $n = $values.Count
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $ws.Range($rangeAddr).Value2 = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: extend the table with two more rows (30 and 31), copying the
# formatting (bordered/bold index cell in column A) from the last existing
# data row so no new style entries are minted. -------------------------------
$ws.Range("A29:W29").Copy($ws.Range("A30:W30"))
$ws.Range("A29:W29").Copy($ws.Range("A31:W31"))

# --- Step 2: column A index values are sequential (row number - 2) for
# every data row from 4 to 31. ------------------------------------------------
for ($r = 4; $r -le 31; $r++) {
    $ws.Range("A" + $r).Value2 = $r - 2
}

# --- Step 3: set the category label and the 21 numeric results for every
# data row (4..31) to their final values. -----------------------------------
$ws.Range("B4").Value2 = "Holden"
Set-RowValues $ws "C4:W4" @(0.9968386538464037, 1.005524107968805, 1.009328248411235, 1.009328248411235, 0.9904312671449856, 1.003531610213833, 0.9852032015031532, 0.9972337030717157, 1.009328248411235, 1.004000857888902, 0.9852032015031532, 1.009328248411235, 1.004000857888902, 0.9946020296960276, 1.000419755867653, 0.9995107692677635, 0.9953475710794862, 0.9995107692677635, 0.9988427404124235, 1.000939842012186, 0.9990114562561292)
$ws.Range("B5").Value2 = "Rizzie Spiral"
Set-RowValues $ws "C5:W5" @(0.9856508305758456, 1.024814008464684, 1.044318817442543, 1.044318817442543, 0.9565789207091203, 1.016633551037466, 0.9331859467236912, 0.9867177960707471, 1.044318817442543, 1.016829456448513, 0.9331859467236912, 1.044318817442543, 1.016829456448513, 0.9750077015861019, 1.001240143512179, 0.9981114068715823, 0.9785554112493499, 0.9981114068715823, 0.9949962627976481, 1.004860773726627, 0.9955911659340761)
$ws.Range("B6").Value2 = "RotRing OmegaMax-90"
Set-RowValues $ws "C6:W6" @(0.9944567125369514, 1.009572868225955, 1.017221081942067, 1.017221081942067, 0.9832263725247519, 1.006456320847752, 0.9742063800023989, 0.9948321550494886, 1.017221081942067, 1.006434153933914, 0.9742063800023989, 1.017221081942067, 1.006434153933914, 0.9903202669681566, 1.000445433235433, 0.9992872052927932, 0.9916990821577549, 0.9992872052927932, 0.9980795821038329, 1.00190788207148, 0.9983007556329099)
$ws.Range("B7").Value2 = "Equal Angle"
Set-RowValues $ws "C7:W7" @(0.9959171899135427, 1.007115500194525, 1.012190248328529, 1.012190248328529, 0.987642972982709, 1.004604643681553, 0.9809153944308359, 0.996374895482707, 1.012190248328529, 1.005070861844384, 0.9809153944308359, 1.012190248328529, 1.005070861844384, 0.9929931281376099, 1.000494025878963, 0.9993921682012497, 0.9939678153962541, 0.9993921682012497, 0.9985234236293229, 1.001256788569164, 0.9987289633573482)
$ws.Range("B8").Value2 = "Tilt Rotate"
Set-RowValues $ws "C8:W8" @(0.9867907332494996, 1.023129713658028, 1.038611432706954, 1.038611432706954, 0.9600163027225374, 1.014644853222266, 0.9381092926669972, 0.988575511280151, 1.038611432706954, 1.016962744990002, 0.9381092926669972, 1.038611432706954, 1.016962744990002, 0.9775360188284998, 1.001876739119751, 0.9978944901213179, 0.9806209236354997, 0.9978944901213179, 0.9951185509033633, 1.003817127264081, 0.9958550730620543)
$ws.Range("B9").Value2 = "CLR"
Set-RowValues $ws "C9:W9" @(0.9995519959253087, 1.000763454272709, 1.001469672136492, 1.001469672136492, 0.9986447945778206, 1.000545559598994, 0.9979290667351292, 0.9995537566894522, 1.001469672136492, 1.000467643747605, 0.9979290667351292, 1.001469672136492, 1.000467643747605, 0.9991983552413672, 1.000009819836457, 0.9999554608730755, 0.9993162354693478, 0.9999554608730755, 0.9998545946361338, 1.000177610136205, 0.9998657429604388)
$ws.Range("B10").Value2 = "Rizzie Hex"
Set-RowValues $ws "C10:W10" @(0.9999670447001311, 1.000057804151613, 1.000095585640128, 1.000095585640128, 0.9999002420402409, 1.000036308079828, 0.9998454599710284, 0.999971771730028, 1.000095585640128, 1.000042821234298, 0.9998454599710284, 1.000095585640128, 1.000042821234298, 0.999944140602663, 1.000004932967214, 0.999994622281818, 0.999951775301819, 0.9999946222818181, 0.9999877278863963, 1.000009299437143, 0.9999896296934119)
$ws.Range("B11").Value2 = "Matthies Hex"
Set-RowValues $ws "C11:W11" @(0.9992874982002709, 1.001217588746727, 1.002311449309538, 1.002311449309538, 0.9978445532733927, 1.000859748084961, 0.9967018547998792, 0.9992998045701011, 1.002311449309538, 1.000761151139993, 0.9967018547998792, 1.002311449309538, 1.000761151139993, 0.9987315029699361, 1.000024324670132, 0.9999248184164699, 0.998916834713381, 0.9999248184164699, 0.9997654883624201, 1.000274680551844, 0.9997854560156079)
$ws.Range("B12").Value2 = "Tilt Rotate_Partial"
Set-RowValues $ws "C12:W12" @(0.9865922006843681, 1.023482383260828, 1.039153382634587, 1.039153382634587, 0.9594151461576711, 1.014853252546956, 0.9371723439895621, 0.9884178872281654, 1.039153382634587, 1.017243490327882, 0.9371723439895621, 1.039153382634587, 1.017243490327882, 0.9772079171587222, 1.001917845506125, 0.9978564056506771, 0.9803360116672709, 0.9978564056506771, 0.9950403544090999, 1.003862960054197, 0.9957912608537525)
$ws.Range("B13").Value2 = "RotRing OmegaMax-60"
Set-RowValues $ws "C13:W13" @(0.9972538424526401, 1.004708238466542, 1.008791869277584, 1.008791869277584, 0.9916917363504226, 1.003277974977679, 0.9872675742589541, 0.9973442244845649, 1.008791869277584, 1.003012341781056, 0.9872675742589541, 1.008791869277584, 1.003012341781056, 0.9951399580200051, 1.000133092116848, 0.9996905951058649, 0.9958445861642168, 0.9996905951058649, 0.9990814069425588, 1.001023499409564, 0.9991684752561804)
$ws.Range("B14").Value2 = "Equal Angle_Partial"
Set-RowValues $ws "C14:W14" @(0.9959168649368393, 1.007152400663158, 1.011914317978948, 1.011914317978948, 0.9876404707368425, 1.00452050773684, 0.980865202526313, 0.9964762379894708, 1.011914317978948, 1.005257440757897, 0.980865202526313, 1.011914317978948, 1.005257440757897, 0.9930613216421051, 1.000587152847368, 0.9993456537543862, 0.9940131694070166, 0.9993456537543861, 0.9984884565499994, 1.001173628835789, 0.9987179304157886)
$ws.Range("B15").Value2 = "Rizzie Hex_Partial"
Set-RowValues $ws "C15:W15" @(1.002613373225497, 0.9954736217128148, 0.9919822668626431, 0.9919822668626431, 1.007908453878432, 0.9969870352280384, 1.012178133907131, 1.002399259494436, 0.9919822668626431, 0.9968986406736717, 1.012178133907131, 0.9919822668626431, 0.9968986406736717, 1.004538387290402, 0.9997560069495844, 1.000353013814482, 1.003896715935434, 1.000353013814482, 1.000918103667236, 0.9991309363063173, 1.000805098122833)
$ws.Range("B16").Value2 = "ND Single"
Set-RowValues $ws "C16:W16" @(0.9768293699999995, 1.0405934, 1.067568100000001, 1.067568100000001, 0.92986294, 1.025639699999999, 0.8914080199999997, 0.980019159999999, 1.067568100000001, 1.029862899999999, 0.8914080199999997, 1.067568100000001, 1.029862899999999, 0.9606354599999993, 1.003346134999999, 0.9962796733333333, 0.9660334299999994, 0.9962796733333331, 0.9914170974999997, 1.006647298, 0.9927229487499998)
$ws.Range("B17").Value2 = "RD Single"
Set-RowValues $ws "C17:W17" @(0.96028664, 1.0671107, 1.1345903, 1.1345903, 0.87989147, 1.049677, 0.8171799800000001, 0.9588603, 1.1345903, 1.0385555, 0.8171799800000001, 1.1345903, 1.0385555, 0.9278677400000001, 0.9994210699999999, 0.9967752600000001, 0.93867404, 0.9967752600000001, 0.987653105, 1.017040544, 0.98826898625)
$ws.Range("B18").Value2 = "TD Single"
Set-RowValues $ws "C18:W18" @(0.96110507, 1.0670152, 1.1220036, 1.1220036, 0.88231288, 1.0456585, 0.81922268, 0.96330981, 1.1220036, 1.0443588, 0.81922268, 1.1220036, 1.0443588, 0.93179074, 1.002731935, 0.9951950266666666, 0.9415621833333333, 0.9951950266666666, 0.9866725374999999, 1.01373875, 0.9881233174999999)
$ws.Range("B19").Value2 = "Morris Single"
Set-RowValues $ws "C19:W19" @(0.98875996, 1.0193501, 1.0353811, 1.0353811, 0.9659909, 1.0132324, 0.94777989, 0.98935163, 1.0353811, 1.0127357, 0.94777989, 1.0353811, 1.0127357, 0.980257795, 1.00074783, 0.99863223, 0.98309185, 0.99863223, 0.9961641625, 1.00400755, 0.99657271)
$ws.Range("B20").Value2 = "Ring Perpendicular to ND"
Set-RowValues $ws "C20:W20" @(0.9910630621917811, 1.015243086849315, 1.029214869863014, 1.029214869863014, 0.9729653639726028, 1.010851709589041, 0.9586703263013698, 0.9911358361643837, 1.029214869863014, 1.009397679726028, 0.9586703263013698, 1.029214869863014, 1.009397679726028, 0.9840340030136987, 1.000230370958904, 0.9990942919634705, 0.9863770227397262, 0.9990942919634705, 0.9970864845205482, 1.003512161589041, 0.9973177418321919)
$ws.Range("B21").Value2 = "Ring Perpendicular to RD"
Set-RowValues $ws "C21:W21" @(0.9838101831578948, 1.02841965, 1.046783046842105, 1.046783046842105, 0.950991384736842, 1.01778432368421, 0.9240491821052631, 0.9861961394736843, 1.046783046842105, 1.021153789473684, 0.9240491821052631, 1.046783046842105, 1.021153789473684, 0.9726014857894736, 1.002481986315789, 0.9973286728070173, 0.976337718245614, 0.9973286728070175, 0.9939490503947368, 1.00451584968421, 0.9948984624342105)
$ws.Range("B22").Value2 = "Ring Perpendicular to TD"
Set-RowValues $ws "C22:W22" @(0.9844748889473686, 1.026846393157895, 1.04795975368421, 1.04795975368421, 0.9530205294736841, 1.017999425789474, 0.9277119763157895, 0.985626017368421, 1.04795975368421, 1.01820266, 0.9277119763157895, 1.04795975368421, 1.01820266, 0.9729573181578947, 1.001338774473684, 0.9979581299999999, 0.9767965084210526, 0.9979581299999999, 0.9945873197368421, 1.005261806526316, 0.9952302055921052)
$ws.Range("B23").Value2 = "OffsetFTD"
Set-RowValues $ws "C23:W23" @(0.9993777859702941, 1.001752986563079, 0.9967624633574996, 0.9967624633574996, 0.998088828514448, 0.9991467770460045, 0.996195832254395, 1.001317684068608, 0.9967624633574996, 1.004198720700661, 0.996195832254395, 0.9967624633574996, 1.004198720700661, 1.000197276477528, 1.001788253335478, 0.9990523387708521, 0.9999241129751169, 0.9990523387708521, 0.9991337005707126, 0.9986594531280699, 0.9996051348093737)
$ws.Range("B24").Value2 = "OffsetATD"
Set-RowValues $ws "C24:W24" @(1.000035746338731, 0.9996602188557281, 1.002007939115971, 1.002007939115971, 1.000119795477235, 1.000605025881231, 1.000538818749538, 0.9992555882860785, 1.002007939115971, 0.9985337662525966, 1.000538818749538, 1.002007939115971, 0.9985337662525966, 0.9995362925010673, 0.9992847562956639, 1.000360174706035, 0.999702777113622, 1.000360174706035, 1.000279067614209, 1.000624841914562, 1.000094612369639)
$ws.Range("B25").Value2 = "OffsetF45"
Set-RowValues $ws "C25:W25" @(1.005818068488067, 0.9895609093610205, 0.9849104901507927, 0.9849104901507927, 1.01762146140804, 0.9941346584451078, 1.027596943453301, 1.004328343866613, 0.9849104901507927, 0.9912397611167881, 1.027596943453301, 0.9849104901507927, 0.9912397611167881, 1.009418352285045, 0.9985289148024277, 1.001249064906961, 1.008218257686052, 1.001249064906961, 1.002391315802238, 0.9988951506719486, 1.001901329536216)
$ws.Range("B26").Value2 = "OffsetA45"
Set-RowValues $ws "C26:W26" @(1.001386664406511, 0.9977681729269999, 0.9944510362848364, 0.9944510362848364, 1.004189154967502, 0.9980061916441396, 1.006234178046038, 1.001748261793105, 0.9944510362848364, 0.999224934269059, 1.006234178046038, 0.9944510362848364, 0.999224934269059, 1.002729556157549, 1.000305799337785, 0.9999700495333114, 1.002281925573869, 0.9999700495333111, 1.000324203251611, 0.9991495698582561, 1.000376074292274)
$ws.Range("B27").Value2 = "OffsetFRD"
Set-RowValues $ws "C27:W27" @(1.001235954354288, 0.9974409140473626, 0.9993968025386902, 0.9993968025386902, 1.003757681824231, 0.9995481699713024, 1.006319995567821, 0.9999643600709193, 0.9993968025386902, 0.9963893061932103, 1.006319995567821, 0.9993968025386902, 0.9963893061932103, 1.001354650880516, 0.9988126302737493, 1.000702034766574, 1.001315085371773, 1.000702034766574, 1.000835514663502, 1.00054777223854, 1.000506648070978)
$ws.Range("B28").Value2 = "OffsetARD"
Set-RowValues $ws "C28:W28" @(1.000188620586372, 0.9998483544216437, 0.998087289624052, 0.998087289624052, 1.000563461108115, 0.9993754327668292, 1.000644452838345, 1.000662801811344, 0.998087289624052, 1.000673115966052, 1.000644452838345, 0.998087289624052, 1.000673115966052, 1.000658784402198, 1.000430868276212, 0.9998016194761495, 1.000502063130256, 0.9998016194761497, 0.9998983697537052, 0.9995361537277745, 1.000005441140344)
$ws.Range("B29").Value2 = "Gaussian Quadrature"
Set-RowValues $ws "C29:W29" @(0.9986405543001359, 1.002335496753226, 1.00431627673697, 1.00431627673697, 0.9958868630628936, 1.001611731054661, 0.9936905877891872, 0.9986985232780632, 1.00431627673697, 1.001515473139585, 0.9936905877891872, 1.00431627673697, 1.001515473139585, 0.9976030304643861, 1.00007801371986, 0.9998407792219141, 0.9979488717429694, 0.9998407792219141, 0.9995407229914696, 1.00049583374057, 0.9995869382643402)
$ws.Range("B30").Value2 = "Michael-CCHex"
Set-RowValues $ws "C30:W30" @(0.9983175073790799, 1.002890120456655, 1.005344394480356, 1.005344394480356, 0.9949095512680121, 1.001995462239654, 0.9921918324376, 0.9983883455549354, 1.005344394480356, 1.001873887087524, 0.9921918324376, 1.005344394480356, 1.001873887087524, 0.9970328597625621, 1.000095697233302, 0.99980337133516, 0.9974610756347347, 0.99980337133516, 0.9994319053461399, 1.000614403172983, 0.9994888876129769)
$ws.Range("B31").Value2 = "Michael-SNHex"
Set-RowValues $ws "C31:W31" @(1.002573931768366, 0.9953952171782479, 0.9932212693825084, 0.9932212693825084, 1.00779525092572, 0.9973736966612945, 1.012190903277881, 1.001952711078432, 0.9932212693825084, 0.9961937151710825, 1.012190903277881, 0.9932212693825084, 0.9961937151710825, 1.004192309224482, 0.9993838234697242, 1.000535295943824, 1.003652850072443, 1.000535295943824, 1.001044954899959, 0.9994802177964692, 1.000837086930442)

Write-Output "done"
